$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates: force text storage so numeric-looking
# strings (e.g. "242.07", "0.0000191") are not coerced into Excel numbers,
# matching the inline-string cells in the source workbook. ---
$priceUpdates = @{
    'D2' = '98.035.72'
    'D3' = '3.146.04'
    'D5' = '242.07'
    'D6' = '609.24'
    'D8' = '0.382'
    'D10' = '3.145.27'
    'D11' = '0.788'
    'D13' = '97.592.90'
    'D15' = '33.94'
    'D17' = '3.729.84'
    'D18' = '3.142.53'
    'D19' = '522.22'
    'D20' = '3.43'
    'D21' = '14.51'
    'D22' = '5.66'
    'D23' = '0.0000191'
    'D24' = '8.76'
    'D25' = '88.73'
    'D26' = '5.44'
    'D27' = '11.57'
    'D28' = '3.314.70'
    'D31' = '0.175'
    'D33' = '0.999'
    'D34' = '8.92'
    'D38' = '24.35'
    'D40' = '0.434'
    'D41' = '466.20'
    'D45' = '3.11'
    'D46' = '163.15'
    'D47' = '1.91'
    'D48' = '0.695'
    'D49' = '4.47'
    'D50' = '44.03'
    'D51' = '0.998'
}
foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = $priceUpdates[$addr]
}
foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Style = "Normal"
}

# --- Volume(1h) (column E) updates ---
$volumeUpdates = @{
    'E2' = '  +5.11%  '
    'E3' = '  +1.46%  '
    'E4' = '  +0.02%  '
    'E5' = '  +2.19%  '
    'E6' = '  -0.59%  '
    'E7' = '  -0.99%  '
    'E8' = '  -1.69%  '
    'E9' = '  +0.11%  '
    'E10' = '  +1.56%  '
    'E11' = '  -5.14%  '
    'E12' = '  +0.21%  '
    'E13' = '  +4.90%  '
    'E14' = '  -1.61%  '
    'E15' = '  -2.85%  '
    'E16' = '  -0.03%  '
    'E17' = '  +1.48%  '
    'E18' = '  +1.21%  '
    'E19' = '  +18.37%  '
    'E20' = '  -7.07%  '
    'E21' = '  -1.19%  '
    'E22' = '  -5.31%  '
    'E23' = '  -3.92%  '
    'E24' = '  -2.86%  '
    'E25' = '  +3.45%  '
    'E26' = '  -4.20%  '
    'E27' = '  -9.39%  '
    'E28' = '  +1.39%  '
    'E29' = '  +0.22%  '
    'E30' = '  -4.48%  '
    'E31' = '  -3.91%  '
    'E32' = '  -1.28%  '
    'E33' = '  -0.81%  '
    'E34' = '  -2.88%  '
    'E35' = '  +2.94%  '
    'E36' = '  -5.06%  '
    'E37' = '  -9.29%  '
    'E38' = '  +1.57%  '
    'E39' = '  -1.21%  '
    'E40' = '  -2.98%  '
    'E41' = '  -2.11%  '
    'E42' = '  -5.85%  '
    'E43' = '  -11.66%  '
    'E44' = '  +0.00%  '
    'E45' = '  -5.02%  '
    'E46' = '  +2.65%  '
    'E47' = '  +2.85%  '
    'E48' = '  -0.70%  '
    'E49' = '  +2.30%  '
    'E50' = '  +0.43%  '
    'E51' = '  +0.00%  '
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}

# --- Row 51: Mantle -> FirstDigitalUSD (coin name + link changed) ---
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
